$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 46
$dateCell = $ws.Cells.Item($row, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = "01/09/2026"
$dateCell.ClearFormats()
$ws.Cells.Item($row, 2).Value = 12871.43
$ws.Cells.Item($row, 3).Value = 0.2097603367030337
$ws.Cells.Item($row, 4).Value = 0.7902396632969663
$ws.Cells.Item($row, 5).Value = -119.71
$ws.Cells.Item($row, 6).Value = -18.85
$ws.Cells.Item($row, 7).Value = -20520.24
$ws.Cells.Item($row, 8).Value = -66.86
$ws.Cells.Item($row, 9).Value = -402.93
$ws.Cells.Item($row, 10).Value = -12.99
